$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.187.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.521.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.39%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.519.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0996"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.333"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.967.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.148.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("E17").Value = "  +3.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.522.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.34%  "
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.88%  "
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("E28").Value = "  +3.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0767"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.92%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.85%  "
$ws.Range("E32").Value = "  +3.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.59%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("E39").Value = "  +4.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.793"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "278.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "132.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.595"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0934"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.86%  "
$ws.Range("E48").Value = "  +5.61%  "
$ws.Range("E49").Value = "  +4.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.759.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.22%  "
